$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A48").Value = "Team Pesto nel Tigullio"
$ws.Range("B48").Value = "MATTEO PILATI | Pinguini Trentini"
$ws.Range("C48").Value = "Federico  Manica | iMontagna"
$ws.Range("D48").Value = "Alessandro Comper | F.C. Gorillaz"
$ws.Range("E48").Value = "Michele Merighi | Clitoriders"
$ws.Range("F48").Value = "Maverick  Bertolini | A.C. Denti"
